$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 15: "Seq" baseline entries for both the 2D (C/D) and DNA (F/G) tables ---
$ws.Range("C15").Value = "Seq"
$ws.Range("D15").Value = 24.6789929867
$ws.Range("F15").Value = "Seq"
$ws.Range("G15").Value = 20.8970220089

# --- Row 17: second mini-table headers (series titles) ---
$ws.Range("C17").Value = "2D MPI kmeans"
$ws.Range("F17").Value = "DNA MPI kmeans"

# --- Row 18: second mini-table column headers ---
$ws.Range("C18").Value = "No. of processes"
$ws.Range("D18").Value = "Time (s)"
$ws.Range("F18").Value = "No. of processes"
$ws.Range("G18").Value = "Time (s)"

# --- Rows 19-33: second mini-table data (15 points), used by both charts ---
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 14.7617080212
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 11.9622499943

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 7.9341380596200004
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 6.6400032043500001

$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 5.3979828357699997
$ws.Range("F21").Value = 6
$ws.Range("G21").Value = 4.9877169132199999

$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 4.2621190547900003
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 3.93734312057

$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 3.5970408916499998
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 3.4324958324399999

$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 3.19927883148
$ws.Range("F24").Value = 12
$ws.Range("G24").Value = 3.0540101528200001

$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 2.9268012046799998
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 2.8865480422999998

$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 2.6982579231299999
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 2.6704561710400001

$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 2.5080280304000002
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 2.5308120250699999

$ws.Range("C28").Value = 20
$ws.Range("D28").Value = 2.4221029281600002
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 2.4044420719100001

$ws.Range("C29").Value = 22
$ws.Range("D29").Value = 2.39634394646
$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 2.3494160175299998

$ws.Range("C30").Value = 24
$ws.Range("D30").Value = 2.3590221405
$ws.Range("F30").Value = 24
$ws.Range("G30").Value = 2.3045308589900002

$ws.Range("C31").Value = 26
$ws.Range("D31").Value = 2.2840700149500002
$ws.Range("F31").Value = 26
$ws.Range("G31").Value = 2.3047590255700001

$ws.Range("C32").Value = 28
$ws.Range("D32").Value = 2.3021240234399998
$ws.Range("F32").Value = 28
$ws.Range("G32").Value = 2.32808709145

$ws.Range("C33").Value = 30
$ws.Range("D33").Value = 2.3463740348800002
$ws.Range("F33").Value = 30
$ws.Range("G33").Value = 2.3625891208600001

# --- Update the active selection to match the edited-file state ---
$ws.Range("G37").Select()
